$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 19, shifting existing rows 19-36 down to 20-37
$ws.Rows("19:19").Insert()

# Populate the newly inserted row 19 with the new weekly record
$ws.Range("A19").Value = 6
$ws.Range("B19").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C19").Value = "Metropolitana"
$ws.Range("D19").Value = 44789
$ws.Range("E19").Value = 13
$ws.Range("F19").Value = 100112035
$ws.Range("G19").Value = "Bruselas (repollito)"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 400
$ws.Range("K19").Value = 15000
$ws.Range("L19").Value = 16000
$ws.Range("M19").Value = 15425
$ws.Range("N19").Value = "`$/malla 15 kilos"
$ws.Range("O19").Value = "Provincia de Quillota"
$ws.Range("P19").Value = 1028
$ws.Range("Q19").Value = 15
$ws.Range("R19").Value = "Hortaliza"
